$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Rabanito" at Vega Central
# Mapocho de Santiago. Insert a new row above the current row 130, pushing
# the existing rows 130:187 down to 131:188.
$ws.Rows(130).Insert()

# Populate the newly inserted row 130 with the new observation's data.
$ws.Cells.Item(130, 1).Value = 9
$ws.Cells.Item(130, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(130, 3).Value = "Metropolitana"
$ws.Cells.Item(130, 4).Value = 44523
$ws.Cells.Item(130, 5).Value = 13
$ws.Cells.Item(130, 6).Value = 300000001
$ws.Cells.Item(130, 7).Value = "Rabanito"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 6100
$ws.Cells.Item(130, 11).Value = 2500
$ws.Cells.Item(130, 12).Value = 3000
$ws.Cells.Item(130, 13).Value = 2750
$ws.Cells.Item(130, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(130, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(130, 16).Value = 28
$ws.Cells.Item(130, 17).Value = 100
$ws.Cells.Item(130, 18).Value = "Hortaliza"
